$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 235.2675016666667
$ws.Range("H2").Value = 705.802505
$ws.Range("I2").Value = 0.5738994362335403
$ws.Range("J2").Value = 0.5738994362335402
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4910443333333334
$ws.Range("N2").Value = 1.473133
$ws.Range("O2").Value = 0.7844104380534107
$ws.Range("P2").Value = 0.7844104380534107
$ws.Range("Q2").Value = 115.5267735109072
$ws.Range("R2").Value = 1039.740961598165
$ws.Range("S2").Value = 0.4501727081745568
$ws.Range("T2").Value = 0.4501727081745567

# Row 3
$ws.Range("G3").Value = 235.2675016666667
$ws.Range("H3").Value = 705.802505
$ws.Range("I3").Value = 0.5738994362335403
$ws.Range("J3").Value = 0.5738994362335402
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.13496
$ws.Range("N3").Value = 0.40488
$ws.Range("O3").Value = 0.2155895619465893
$ws.Range("P3").Value = 0.2155895619465893
$ws.Range("Q3").Value = 31.75170202493333
$ws.Range("R3").Value = 285.7653182244
$ws.Range("S3").Value = 0.1237267280589835
$ws.Range("T3").Value = 0.1237267280589835

# Row 4
$ws.Range("I4").Value = 0.3286113026040369
$ws.Range("J4").Value = 0.3286113026040369
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4910443333333334
$ws.Range("N4").Value = 1.473133
$ws.Range("O4").Value = 0.7844104380534107
$ws.Range("P4").Value = 0.7844104380534107
$ws.Range("Q4").Value = 66.14992302172622
$ws.Range("R4").Value = 595.349307195536
$ws.Range("S4").Value = 0.2577661358249345
$ws.Range("T4").Value = 0.2577661358249345

# Row 5
$ws.Range("I5").Value = 0.3286113026040369
$ws.Range("J5").Value = 0.3286113026040369
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.13496
$ws.Range("N5").Value = 0.40488
$ws.Range("O5").Value = 0.2155895619465893
$ws.Range("P5").Value = 0.2155895619465893
$ws.Range("Q5").Value = 18.18083013077333
$ws.Range("R5").Value = 163.62747117696
$ws.Range("S5").Value = 0.07084516677910241
$ws.Range("T5").Value = 0.07084516677910241

# Row 6
$ws.Range("G6").Value = 0.325805
$ws.Range("H6").Value = 0.977415
$ws.Range("I6").Value = 0.0007947519504286909
$ws.Range("J6").Value = 0.0007947519504286907
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.4910443333333334
$ws.Range("N6").Value = 1.473133
$ws.Range("O6").Value = 0.7844104380534107
$ws.Range("P6").Value = 0.7844104380534107
$ws.Range("Q6").Value = 0.1599846990216667
$ws.Range("R6").Value = 1.439862291195
$ws.Range("S6").Value = 0.000623411725579572
$ws.Range("T6").Value = 0.0006234117255795718

# Row 7
$ws.Range("G7").Value = 0.325805
$ws.Range("H7").Value = 0.977415
$ws.Range("I7").Value = 0.0007947519504286909
$ws.Range("J7").Value = 0.0007947519504286907
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.13496
$ws.Range("N7").Value = 0.40488
$ws.Range("O7").Value = 0.2155895619465893
$ws.Range("P7").Value = 0.2155895619465893
$ws.Range("Q7").Value = 0.0439706428
$ws.Range("R7").Value = 0.3957357852
$ws.Range("S7").Value = 0.0001713402248491189
$ws.Range("T7").Value = 0.0001713402248491189

# Row 8
$ws.Range("G8").Value = 39.46134166666666
$ws.Range("H8").Value = 118.384025
$ws.Range("I8").Value = 0.09625996610278018
$ws.Range("J8").Value = 0.09625996610278018
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.4910443333333334
$ws.Range("N8").Value = 1.473133
$ws.Range("O8").Value = 0.7844104380534107
$ws.Range("P8").Value = 0.7844104380534107
$ws.Range("Q8").Value = 19.37726821114722
$ws.Range("R8").Value = 174.395413900325
$ws.Range("S8").Value = 0.07550732217768827
$ws.Range("T8").Value = 0.07550732217768827

# Row 9
$ws.Range("G9").Value = 39.46134166666666
$ws.Range("H9").Value = 118.384025
$ws.Range("I9").Value = 0.09625996610278018
$ws.Range("J9").Value = 0.09625996610278018
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.13496
$ws.Range("N9").Value = 0.40488
$ws.Range("O9").Value = 0.2155895619465893
$ws.Range("P9").Value = 0.2155895619465893
$ws.Range("Q9").Value = 5.325702671333333
$ws.Range("R9").Value = 47.931324042
$ws.Range("S9").Value = 0.02075264392509191
$ws.Range("T9").Value = 0.02075264392509191

# Row 10
$ws.Range("G10").Value = 0.178139
$ws.Range("H10").Value = 0.534417
$ws.Range("I10").Value = 0.0004345431092138444
$ws.Range("J10").Value = 0.0004345431092138443
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.4910443333333334
$ws.Range("N10").Value = 1.473133
$ws.Range("O10").Value = 0.7844104380534107
$ws.Range("P10").Value = 0.7844104380534107
$ws.Range("Q10").Value = 0.08747414649566668
$ws.Range("R10").Value = 0.7872673184610001
$ws.Range("S10").Value = 0.0003408601506515228
$ws.Range("T10").Value = 0.0003408601506515227

# Row 11
$ws.Range("G11").Value = 0.178139
$ws.Range("H11").Value = 0.534417
$ws.Range("I11").Value = 0.0004345431092138444
$ws.Range("J11").Value = 0.0004345431092138443
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.13496
$ws.Range("N11").Value = 0.40488
$ws.Range("O11").Value = 0.2155895619465893
$ws.Range("P11").Value = 0.2155895619465893
$ws.Range("Q11").Value = 0.02404163944
$ws.Range("R11").Value = 0.21637475496
$ws.Range("S11").Value = 0.000093682958562321615148503484
$ws.Range("T11").Value = 0.000093682958562321601595976328
